$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsLpWOU = $wb.Worksheets.Item("LpWOU")

# Update the unit note text on the About sheet (cubic km -> billion cubic m)
$wsAbout.Range("A9").Value = "For the U.S., the water output unit is billion cubic m, which is equivalent to Tl (teraliters, or 10^12 liters)"

# Leave a stray selection on the About sheet (matches author's saved state)
[void]$wsAbout.Range("J10").Select()

# Make LpWOU the active/selected sheet, with B2 selected
$wsLpWOU.Activate()
[void]$wsLpWOU.Range("B2").Select()
